$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

Set-TextValue "D2" "303.62"
Set-TextValue "E2" "4.99%"
Set-TextValue "G2" "15"
Set-TextValue "D3" "34.96"
Set-TextValue "E3" "12.68%"
Set-TextValue "G3" "15"
Set-TextValue "E4" "4.51%"
Set-TextValue "G4" "15"
Set-TextValue "D5" "0.07831"
Set-TextValue "E5" "6.29%"
Set-TextValue "G5" "15"
Set-TextValue "D6" "2.322"
Set-TextValue "E6" "0.75%"
Set-TextValue "G6" "15"
Set-TextValue "D7" "8.058"
Set-TextValue "E7" "5.01%"
Set-TextValue "G7" "15"
Set-TextValue "D8" "3.971"
Set-TextValue "E8" "6.28%"
Set-TextValue "G8" "15"
Set-TextValue "D9" "0.9252"
Set-TextValue "E9" "0.81%"
Set-TextValue "G9" "15"
Set-TextValue "D10" "0.1010"
Set-TextValue "E10" "10.14%"
Set-TextValue "G10" "15"
Set-TextValue "D11" "0.1827"
Set-TextValue "E11" "7.29%"
Set-TextValue "G11" "15"
Set-TextValue "D12" "0.08538"
Set-TextValue "E12" "2.75%"
Set-TextValue "G12" "15"
Set-TextValue "D13" "0.03419"
Set-TextValue "E13" "10.21%"
Set-TextValue "G13" "15"
Set-TextValue "D14" "0.09909"
Set-TextValue "E14" "-0.83%"
Set-TextValue "G14" "15"
Set-TextValue "D15" "0.001480"
Set-TextValue "E15" "-1.42%"
Set-TextValue "G15" "15"
Set-TextValue "D16" "0.005822"
Set-TextValue "E16" "1.38%"
Set-TextValue "G16" "15"
Set-TextValue "E17" "0.08%"
Set-TextValue "G17" "15"
Set-TextValue "D18" "2.127"
Set-TextValue "E18" "0.13%"
Set-TextValue "G18" "15"
Set-TextValue "E19" "2.96%"
Set-TextValue "G19" "15"
Set-TextValue "D20" "0.1326"
Set-TextValue "E20" "2.77%"
Set-TextValue "G20" "15"
Set-TextValue "D21" "4.554"
Set-TextValue "E21" "9.67%"
Set-TextValue "G21" "15"
Set-TextValue "D22" "0.2221"
Set-TextValue "E22" "4.68%"
Set-TextValue "G22" "15"
Set-TextValue "D23" "0.04634"
Set-TextValue "E23" "2.98%"
Set-TextValue "G23" "15"
Set-TextValue "D24" "0.001218"
Set-TextValue "E24" "0.25%"
Set-TextValue "G24" "15"
Set-TextValue "D25" "0.004458"
Set-TextValue "E25" "6.21%"
Set-TextValue "G25" "15"
Set-TextValue "G26" "15"
Set-TextValue "D27" "0.0003398"
Set-TextValue "G27" "15"
Set-TextValue "G28" "15"
Set-TextValue "G29" "15"
Set-TextValue "G30" "15"
Set-TextValue "G31" "15"
Set-TextValue "G32" "15"
Set-TextValue "G33" "15"
Set-TextValue "G34" "15"
Set-TextValue "G35" "15"
Set-TextValue "G36" "15"
Set-TextValue "G37" "15"
Set-TextValue "G38" "15"
Set-TextValue "D39" "0.01746"
Set-TextValue "E39" "10.27%"
Set-TextValue "G39" "15"
Set-TextValue "D40" "0.04727"
Set-TextValue "E40" "5.14%"
Set-TextValue "G40" "15"
Set-TextValue "D41" "0.007775"
Set-TextValue "E41" "5.01%"
Set-TextValue "G41" "15"
Set-TextValue "D42" "0.1414"
Set-TextValue "E42" "5.79%"
Set-TextValue "G42" "15"
Set-TextValue "D43" "0.008839"
Set-TextValue "E43" "-10.31%"
Set-TextValue "G43" "15"
Set-TextValue "D44" "0.002290"
Set-TextValue "E44" "2.78%"
Set-TextValue "G44" "15"
Set-TextValue "D45" "0.009971"
Set-TextValue "E45" "17.03%"
Set-TextValue "G45" "15"
Set-TextValue "D46" "0.00006088"
Set-TextValue "E46" "-0.19%"
Set-TextValue "G46" "15"
Set-TextValue "G47" "15"
Set-TextValue "D48" "3.906"
Set-TextValue "E48" "49.90%"
Set-TextValue "G48" "15"
Set-TextValue "D49" "0.002689"
Set-TextValue "G49" "15"
Set-TextValue "G50" "15"
Set-TextValue "G51" "15"
